$d = $word.ActiveDocument

# 1. "Pillajo Magaly" -- collapse the separate " " + "Magaly" runs (with
#    their spell-check proofErr wrapper) into a single run's text " Magaly"
#    by simply re-finding the two words and replacing them as one block.
$d.Content.Find.Execute("Pillajo Magaly", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pillajo Magaly", 2) | Out-Null

# 2. Split the heading "CASOS DE USO DEL PROCESO DE INGRESO DE PROVEEDORES"
#    so that "INGRESO" becomes "REGISTRO", each segment now living in its
#    own run (bold kept throughout).
$rng = $d.Content
$rng.Find.Execute("CASOS DE USO DEL PROCESO DE INGRESO DE PROVEEDORES", $true, $false, `
                   $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$heading = $rng
$heading.Text = "CASOS DE USO DEL PROCESO DE "
$afterStart = $heading.End
$r2 = $d.Range($afterStart, $afterStart)
$r2.InsertAfter("REGISTRO")
$r2.Bold = 1
$r3End = $afterStart + 8
$r3 = $d.Range($r3End, $r3End)
$r3.InsertAfter(" DE PROVEEDORES")
$r3.Bold = 1

# 3. Remove the leftover "_GoBack" bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
